# Update "want-to-go" counts (column F) and a couple of "lowest price"
# values (column G) across the workbook's sheets, per the commit's
# refreshed scrape data.
#
# Sheet order in the workbook:
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Performances)
#   3 = 本地生活 (Local Life)
#   4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 962
$ws1.Range("F6").Value  = 361
$ws1.Range("F9").Value  = 1441
$ws1.Range("F11").Value = 1329
$ws1.Range("F12").Value = 2996
$ws1.Range("F13").Value = 402
$ws1.Range("F14").Value = 1607
$ws1.Range("F16").Value = 789
$ws1.Range("F18").Value = 1374
$ws1.Range("F23").Value = 3464
$ws1.Range("F26").Value = 1531

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value  = 26
$ws2.Range("F7").Value  = 49
$ws2.Range("G12").Value = 224

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 794

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 794
$ws4.Range("F9").Value  = 26
$ws4.Range("F11").Value = 49
$ws4.Range("F15").Value = 962
$ws4.Range("F16").Value = 361
$ws4.Range("F19").Value = 1441
$ws4.Range("F21").Value = 1329
$ws4.Range("F22").Value = 2996
$ws4.Range("F23").Value = 402
$ws4.Range("F24").Value = 1607
$ws4.Range("F26").Value = 789
$ws4.Range("F28").Value = 1374
$ws4.Range("F35").Value = 3464
$ws4.Range("F38").Value = 1531
$ws4.Range("G39").Value = 224
